# Apply coin price/volume/name/link updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "277.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.93%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.30"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.81%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.885"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.03%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06419"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.20%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.965"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.03%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.182"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-5.67%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8832"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.78%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1543"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.88%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05139"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.52%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07425"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.31%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02886"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-4.53%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08976"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.55%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001566"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.42%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006387"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.25%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006166"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "3.76%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.477"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.78%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.315"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.04%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.10%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.48%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.901"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.28%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04413"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.95%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.001176"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.48%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.003866"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-9.04%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001180"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.63%"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "15.61%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04140"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.86%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006785"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.83%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.52%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001911"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-10.88%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01146"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "6.97%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005321"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.00%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.687"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "13.27%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-7.32%"
